# Re-label the survey header row (row 1).
#
# The sheet has one header row (row 1) with survey question codes in
# A1:AK1, backed by the shared-string table, followed by 19 rows of
# response data (rows 2-20) that are untouched by this change.
#
# Original headers: v1, v2, v3, ..., v37   (A1:AK1)
# New headers:       q1, q3, q4a, q5a, q7pa, q7sa, q8, q10, v9, v10, ..., v37
#
# i.e. the old v1..v8 columns are replaced by eight new "q*" columns,
# and v9..v37 shift left to directly follow them (I1:AK1 keeps the same
# v9..v37 labels the old I1:AK1 had). Only setting the header cells'
# values is required: doing so automatically drops the now-unused
# v1..v8 strings from the workbook's shared-string table and appends
# the new q* strings, while leaving all the underlying response data in
# rows 2-20 completely unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "q1"
$ws.Range("B1").Value = "q3"
$ws.Range("C1").Value = "q4a"
$ws.Range("D1").Value = "q5a"
$ws.Range("E1").Value = "q7pa"
$ws.Range("F1").Value = "q7sa"
$ws.Range("G1").Value = "q8"
$ws.Range("H1").Value = "q10"

$ws.Range("I1").Value = "v9"
$ws.Range("J1").Value = "v10"
$ws.Range("K1").Value = "v11"
$ws.Range("L1").Value = "v12"
$ws.Range("M1").Value = "v13"
$ws.Range("N1").Value = "v14"
$ws.Range("O1").Value = "v15"
$ws.Range("P1").Value = "v16"
$ws.Range("Q1").Value = "v17"
$ws.Range("R1").Value = "v18"
$ws.Range("S1").Value = "v19"
$ws.Range("T1").Value = "v20"
$ws.Range("U1").Value = "v21"
$ws.Range("V1").Value = "v22"
$ws.Range("W1").Value = "v23"
$ws.Range("X1").Value = "v24"
$ws.Range("Y1").Value = "v25"
$ws.Range("Z1").Value = "v26"
$ws.Range("AA1").Value = "v27"
$ws.Range("AB1").Value = "v28"
$ws.Range("AC1").Value = "v29"
$ws.Range("AD1").Value = "v30"
$ws.Range("AE1").Value = "v31"
$ws.Range("AF1").Value = "v32"
$ws.Range("AG1").Value = "v33"
$ws.Range("AH1").Value = "v34"
$ws.Range("AI1").Value = "v35"
$ws.Range("AJ1").Value = "v36"
$ws.Range("AK1").Value = "v37"

# Match the author's final cursor position/selection (J24) recorded in
# the saved sheet view.
$ws.Range("J24").Select()
